$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "Timestamp" -> "Transaction Date"
$ws.Range("A1").Value = "Transaction Date"

# Fix existing row 4, column K text (drop the extra trailing blank line)
$ws.Range("K4").Value = "OTP for transaction amount of Rs1000 is 123456.`n"
$ws.Rows(4).EntireRow.AutoFit()

# New row 5 data
$ws.Range("A5").Value = "2025-08-06 09:48:30"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "5555"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "7777"

$ws.Range("D5").Value = "sirius"
$ws.Range("E5").Value = "Fitness Inspection Renewal - xxxx"

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "990"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "1000"

$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "123456"

$ws.Range("I5").Value = "kreacher"
$ws.Range("J5").Value = "1987d9a587a4073d"
$ws.Range("K5").Value = "OTP for transaction amount of Rs1000 is 123456.`n"

$ws.Rows(5).EntireRow.AutoFit()
